$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card3")

# Remove trailing space from the header text in M1
$ws.Range("M1").Value = "Serviced by"

# Fill in the "Serviced by" column for the data rows
$ws.Range("M2").Value = "nan"
$ws.Range("M3").Value = "nan"
$ws.Range("M4").Value = "nan"
$ws.Range("M5").Value = "nan"
$ws.Range("M6").Value = "nan"
$ws.Range("M7").Value = "nan"
$ws.Range("M8").Value = "م.محمد عبدالله "
$ws.Range("M9").Value = "nan"
$ws.Range("M10").Value = "nan"
$ws.Range("M11").Value = "nan"
$ws.Range("M12").Value = "nan"
$ws.Range("M13").Value = "nan"
